$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for the two new columns ---
$ws.Range("S1").Value = "IF_val"
$ws.Range("T1").Value = "IF_chg"

# --- Fill the bulk (default) data for rows 2..35: zeros formatted with 2 decimals ---
$ws.Range("S2:T35").Value = 0
$ws.Range("S2:T35").NumberFormat = "0.00"

# --- Special rows that carry real (non-zero) IF_val / IF_chg figures ---
# IF_val column (S) uses a left-aligned 2-decimal number format.
$ws.Range("S7").Value = 0.04
$ws.Range("S7").NumberFormat = "0.00"
$ws.Range("S7").HorizontalAlignment = -4131

$ws.Range("S10").Value = -0.15

$ws.Range("S14").Value = 0.02
$ws.Range("S14").NumberFormat = "0.00"
$ws.Range("S14").HorizontalAlignment = -4131

$ws.Range("S27").Value = 0.02
$ws.Range("S27").NumberFormat = "0.00"
$ws.Range("S27").HorizontalAlignment = -4131

$ws.Range("S30").Value = -0.01
$ws.Range("S30").NumberFormat = "0.00"
$ws.Range("S30").HorizontalAlignment = -4131

# IF_chg column (T) uses the plain default/general number format (no explicit style),
# so drop the style picked up from the bulk fill above before writing the real value.
$ws.Range("T7").Style = "Normal"
$ws.Range("T7").Value = -0.27999999999999997

$ws.Range("T10").Style = "Normal"
$ws.Range("T10").Value = 0.079999999999999988

$ws.Range("T14").Style = "Normal"
$ws.Range("T14").Value = -0.06

$ws.Range("T27").Style = "Normal"
$ws.Range("T27").Value = -0.16999999999999998

$ws.Range("T30").Style = "Normal"
$ws.Range("T30").Value = 0.25

# --- New trailing blank (formatted only) row 36 ---
$ws.Range("S36").NumberFormat = "0.00"
$ws.Range("T36").NumberFormat = "0.00"

# --- Update the selection to mirror the authored view state ---
$ws.Range("V21").Select()
